$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated Diebold-Mariano statistics and p-values
$ws.Range("C2").Value = 0.6191354630309253
$ws.Range("D2").Value = 0.5421860390687239

$ws.Range("C3").Value = 1.329899267606514
$ws.Range("D3").Value = 0.1971758856784527

$ws.Range("C4").Value = 0.9459950934441
$ws.Range("D4").Value = 0.3544263108606702

$ws.Range("C5").Value = 1.969104919576889
$ws.Range("D5").Value = 0.06166382251492286
$ws.Range("G5").Value = "No"

$ws.Range("C6").Value = 0.7281727061528513
$ws.Range("D6").Value = 0.4741877812233901

$ws.Range("C7").Value = 0.5230055931100051
$ws.Range("D7").Value = 0.6061983706090444

$ws.Range("C8").Value = 1.827157562554717
$ws.Range("D8").Value = 0.08127450029725947

$ws.Range("C9").Value = -0.3123449568785451
$ws.Range("D9").Value = 0.7577196605340986

$ws.Range("C10").Value = 0.5525072448575221
$ws.Range("D10").Value = 0.5861711684781947

$ws.Range("C11").Value = 0.8883016465124707
$ws.Range("D11").Value = 0.3839880547286343
